# The commit swaps the data held in row 2 and row 3 of the sheet
# (the header row, row 1, is left untouched). Implement this as a
# three-way shuffle through a scratch row so that Excel's "Copy"
# semantics (which, unlike a plain Value swap, keep literal text such
# as the "YYYY-MM-DD" date strings from being reinterpreted as real
# dates) are preserved for every column, A through AY.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = $ws.Range("A2:AY2")
$row3 = $ws.Range("A3:AY3")
$scratch = $ws.Range("A10:AY10")

$scratch.ClearContents()

# scratch <- row2 (original)
$row2.Copy($scratch)
$row2.ClearContents()

# row2 <- row3 (original)
$row3.Copy($row2)
$row3.ClearContents()

# row3 <- scratch (original row2)
$scratch.Copy($row3)
$scratch.ClearContents()
